# Insert a new data row at row 12 (pushing existing rows 12-35 down to 13-36)
# by copying each row's full contents (values + styles) down by one row,
# working from the bottom up so we never overwrite data before reading it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 35; $r -ge 12; $r--) {
    $src = $ws.Range("A" + $r + ":R" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":R" + ($r + 1))
    $src.Copy($dst)
}

# Populate the newly freed row 12 with the new record.
$ws.Cells.Item(12, 1).Value = 7
$ws.Cells.Item(12, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(12, 3).Value = "Ñuble"
$ws.Cells.Item(12, 4).Value = 44536
$ws.Cells.Item(12, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(12, 5).Value = 16
$ws.Cells.Item(12, 6).Value = 100112026
$ws.Cells.Item(12, 7).Value = "Haba"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 80
$ws.Cells.Item(12, 11).Value = 6500
$ws.Cells.Item(12, 12).Value = 7000
$ws.Cells.Item(12, 13).Value = 6750
$ws.Cells.Item(12, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(12, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(12, 16).Value = 270
$ws.Cells.Item(12, 17).Value = 25
$ws.Cells.Item(12, 18).Value = "Hortaliza"
